$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated transition-probability matrix values (recomputed after adding more
# simulated games / sped-up simulate-game logic per the commit message).
$updates = @{
    "B2" = 0.1886792452830189
    "C2" = 0.5440251572327044
    "J2" = 0.009433962264150943
    "O2" = 0.003144654088050315
    "P2" = 0.119496855345912
    "S2" = 0.1352201257861635
    "B3" = 0.02259887005649718
    "C3" = 0.03389830508474576
    "J3" = 0.01694915254237288
    "P3" = 0.7175141242937854
    "S3" = 0.2090395480225989
    "P4" = 0.6444444444444445
    "S4" = 0.3555555555555556
    "B6" = 0.0593607305936073
    "D6" = 0.0319634703196347
    "F6" = 0.0684931506849315
    "J6" = 0.2420091324200913
    "O6" = 0.0136986301369863
    "Q6" = 0.1917808219178082
    "R6" = 0.0547945205479452
    "S6" = 0.3378995433789954
    "B7" = 0.1222222222222222
    "D7" = 0.01666666666666667
    "F7" = 0.05
    "J7" = 0.1555555555555556
    "O7" = 0.02222222222222222
    "Q7" = 0.1277777777777778
    "R7" = 0.05555555555555555
    "S7" = 0.45
    "B8" = 0.1026894865525672
    "D8" = 0.01711491442542787
    "F8" = 0.05623471882640587
    "J8" = 0.117359413202934
    "O8" = 0.02689486552567237
    "Q8" = 0.2200488997555012
    "R8" = 0.05378973105134474
    "S8" = 0.4058679706601467
    "B9" = 0.1019108280254777
    "D9" = 0.01273885350318471
    "F9" = 0.05732484076433121
    "J9" = 0.1146496815286624
    "O9" = 0.03821656050955414
    "Q9" = 0.1719745222929936
    "R9" = 0.07006369426751592
    "S9" = 0.4331210191082803
    "B10" = 0.1271860095389507
    "D10" = 0.02066772655007949
    "E10" = 0.001589825119236884
    "F10" = 0.06756756756756757
    "J10" = 0.1216216216216216
    "O10" = 0.01271860095389507
    "Q10" = 0.2241653418124006
    "R10" = 0.05564387917329094
    "S10" = 0.3688394276629571
    "F11" = 0.009900990099009901
    "G11" = 0.1353135313531353
    "J11" = 0.09900990099009901
    "K11" = 0.2079207920792079
    "L11" = 0.5313531353135313
    "S11" = 0.0165016501650165
    "G12" = 0.7090909090909091
    "J12" = 0.2242424242424242
    "K12" = 0.02424242424242424
    "L12" = 0.02424242424242424
    "S12" = 0.01818181818181818
    "G13" = 0.5510204081632653
    "J13" = 0.3469387755102041
    "S13" = 0.1020408163265306
    "F15" = 0.05116279069767442
    "H15" = 0.2186046511627907
    "I15" = 0.05116279069767442
    "J15" = 0.3255813953488372
    "K15" = 0.05581395348837209
    "M15" = 0.0186046511627907
    "O15" = 0.05116279069767442
    "S15" = 0.2279069767441861
    "F16" = 0.005235602094240838
    "H16" = 0.1727748691099476
    "I16" = 0.07329842931937172
    "J16" = 0.3455497382198953
    "K16" = 0.1256544502617801
    "M16" = 0.03141361256544502
    "O16" = 0.05235602094240838
    "S16" = 0.193717277486911
    "F17" = 0.01948051948051948
    "H17" = 0.1601731601731602
    "I17" = 0.06493506493506493
    "J17" = 0.4632034632034632
    "K17" = 0.08874458874458875
    "M17" = 0.02164502164502164
    "N17" = 0.002164502164502165
    "O17" = 0.06277056277056277
    "F18" = 0.01612903225806452
    "H18" = 0.1612903225806452
    "I18" = 0.0564516129032258
    "J18" = 0.5
    "K18" = 0.09677419354838709
    "M18" = 0.008064516129032258
    "O18" = 0.07258064516129033
    "S18" = 0.08870967741935484
    "F19" = 0.02218570254724733
    "H19" = 0.1922760887428102
    "I19" = 0.07641741988496302
    "J19" = 0.3853738701725555
    "K19" = 0.1199671322925226
    "M19" = 0.02465078060805259
    "N19" = 0.0008216926869350862
    "O19" = 0.07559572719802794
    "S19" = 0.1027115858668858
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

